$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.228.47"
Set-TextValue $ws.Range("E2") "  +1.14%  "
Set-TextValue $ws.Range("D3") "2.085.65"
Set-TextValue $ws.Range("E3") "  -1.24%  "
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  -0.28%  "
Set-TextValue $ws.Range("D5") "341.18"
Set-TextValue $ws.Range("E5") "  -2.04%  "
Set-TextValue $ws.Range("D6") "1.003"
Set-TextValue $ws.Range("E6") "  -0.35%  "
Set-TextValue $ws.Range("D7") "0.5289"
Set-TextValue $ws.Range("E7") "  +1.75%  "
Set-TextValue $ws.Range("D8") "0.4374"
Set-TextValue $ws.Range("E8") "  -2.07%  "
Set-TextValue $ws.Range("D9") "55.00"
Set-TextValue $ws.Range("E9") "  +1.66%  "
Set-TextValue $ws.Range("D10") "0.09361"
Set-TextValue $ws.Range("E10") "  -0.01%  "
Set-TextValue $ws.Range("D11") "1.173"
Set-TextValue $ws.Range("E11") "  -0.73%  "
Set-TextValue $ws.Range("D12") "24.46"
Set-TextValue $ws.Range("E12") "  -2.86%  "
Set-TextValue $ws.Range("D13") "8.503"
Set-TextValue $ws.Range("E13") "  +2.11%  "
Set-TextValue $ws.Range("D14") "6.852"
Set-TextValue $ws.Range("E14") "  +0.09%  "
Set-TextValue $ws.Range("D15") "2.045.65"
Set-TextValue $ws.Range("E15") "  -7.27%  "
Set-TextValue $ws.Range("D16") "101.39"
Set-TextValue $ws.Range("E16") "  -1.28%  "
Set-TextValue $ws.Range("D17") "0.00001158"
Set-TextValue $ws.Range("E17") "  -0.43%  "
Set-TextValue $ws.Range("E18") "  -0.43%  "
Set-TextValue $ws.Range("D19") "21.01"
Set-TextValue $ws.Range("E19") "  -2.12%  "
Set-TextValue $ws.Range("D20") "0.06705"
Set-TextValue $ws.Range("E20") "  +0.43%  "
Set-TextValue $ws.Range("D21") "6.297"
Set-TextValue $ws.Range("E21") "  -0.27%  "
Set-TextValue $ws.Range("D23") "30.296.37"
Set-TextValue $ws.Range("E23") "  +1.24%  "
Set-TextValue $ws.Range("D24") "12.42"
Set-TextValue $ws.Range("E24") "  -2.49%  "
Set-TextValue $ws.Range("D25") "2.319"
Set-TextValue $ws.Range("E25") "  -0.39%  "
Set-TextValue $ws.Range("D26") "21.77"
Set-TextValue $ws.Range("E26") "  -1.78%  "
Set-TextValue $ws.Range("D27") "6.878"
Set-TextValue $ws.Range("E27") "  +7.77%  "
Set-TextValue $ws.Range("D28") "162.10"
Set-TextValue $ws.Range("E28") "  -0.12%  "
Set-TextValue $ws.Range("D29") "2.479"
Set-TextValue $ws.Range("E29") "  -2.85%  "
Set-TextValue $ws.Range("E30") "  -0.30%  "
Set-TextValue $ws.Range("D31") "1.128"
Set-TextValue $ws.Range("E31") "  -1.93%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D32") "0.1047"
Set-TextValue $ws.Range("E32") "  -0.90%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D33") "1.657"
Set-TextValue $ws.Range("E33") "  -7.59%  "
Set-TextValue $ws.Range("D34") "6.241"
Set-TextValue $ws.Range("E34") "  -0.20%  "
Set-TextValue $ws.Range("D35") "3.911"
Set-TextValue $ws.Range("E35") "  -1.11%  "
Set-TextValue $ws.Range("D36") "10.03"
Set-TextValue $ws.Range("E36") "  -7.39%  "
Set-TextValue $ws.Range("D37") "0.02617"
Set-TextValue $ws.Range("E37") "  +0.91%  "
Set-TextValue $ws.Range("D38") "0.06741"
Set-TextValue $ws.Range("E38") "  -0.90%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D39") "0.6960"
Set-TextValue $ws.Range("E39") "  -1.24%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "12.54"
Set-TextValue $ws.Range("E40") "  -1.66%  "
Set-TextValue $ws.Range("D41") "1.337"
Set-TextValue $ws.Range("E41") "  -0.20%  "
Set-TextValue $ws.Range("D42") "0.2199"
Set-TextValue $ws.Range("E42") "  -1.74%  "
Set-TextValue $ws.Range("D43") "0.6724"
Set-TextValue $ws.Range("E43") "  -1.86%  "
Set-TextValue $ws.Range("D44") "2.388"
Set-TextValue $ws.Range("E44") "  +0.97%  "
Set-TextValue $ws.Range("D45") "14.28"
Set-TextValue $ws.Range("E45") "  -1.30%  "
Set-TextValue $ws.Range("D46") "1.002"
Set-TextValue $ws.Range("E46") "  -0.50%  "
Set-TextValue $ws.Range("D47") "1.285"
Set-TextValue $ws.Range("E47") "  +4.99%  "
Set-TextValue $ws.Range("D48") "3.634"
Set-TextValue $ws.Range("E48") "  -0.17%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D49") "1.208"
Set-TextValue $ws.Range("E49") "  +1.72%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D50") "0.00000000341"
Set-TextValue $ws.Range("E50") "  -3.09%  "
Set-TextValue $ws.Range("D51") "1.209"
Set-TextValue $ws.Range("E51") "  -1.31%  "
